$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "2022-Q1" sheet by duplicating "2021-Q4" (same column
#    layout/styling) and dropping it right before "总计" (the last sheet).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($totalSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count - 1)
$newSheet.Name = "2022-Q1"

# Helper: write a value into a cell while forcing a text number-format first
# for values that look numeric (fund codes / decimal strings) so Excel does
# not silently coerce them into real numbers (mirrors the source data, which
# stores these as text).
function Set-TextValue($sheet, $row, $col, $val) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# Row 2 - 006679
Set-TextValue $newSheet 2 2 "006679"
$newSheet.Cells.Item(2, 3).Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇A"
Set-TextValue $newSheet 2 4 "14.75"
Set-TextValue $newSheet 2 5 "83.19"
Set-TextValue $newSheet 2 6 "15.14"
Set-TextValue $newSheet 2 7 "2.2332"
$newSheet.Cells.Item(2, 8).Value = 1

# Row 3 - 162719
Set-TextValue $newSheet 3 2 "162719"
$newSheet.Cells.Item(3, 3).Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A"
Set-TextValue $newSheet 3 4 "14.75"
Set-TextValue $newSheet 3 5 "83.19"
Set-TextValue $newSheet 3 6 "15.14"
Set-TextValue $newSheet 3 7 "2.2332"
$newSheet.Cells.Item(3, 8).Value = 1

# Row 4 - 006680
Set-TextValue $newSheet 4 2 "006680"
$newSheet.Cells.Item(4, 3).Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇C"
Set-TextValue $newSheet 4 4 "4.73"
Set-TextValue $newSheet 4 5 "83.19"
Set-TextValue $newSheet 4 6 "15.14"
Set-TextValue $newSheet 4 7 "0.7161"
$newSheet.Cells.Item(4, 8).Value = 1

# Row 5 - 004243
Set-TextValue $newSheet 5 2 "004243"
$newSheet.Cells.Item(5, 3).Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C"
Set-TextValue $newSheet 5 4 "4.73"
Set-TextValue $newSheet 5 5 "83.19"
Set-TextValue $newSheet 5 6 "15.14"
Set-TextValue $newSheet 5 7 "0.7161"
$newSheet.Cells.Item(5, 8).Value = 1

# Row 6 - 160416
Set-TextValue $newSheet 6 2 "160416"
$newSheet.Cells.Item(6, 3).Value = "华安标普全球石油指数 (QDII-LOF)"
Set-TextValue $newSheet 6 4 "3.37"
Set-TextValue $newSheet 6 5 "95.08"
Set-TextValue $newSheet 6 6 "3.94"
Set-TextValue $newSheet 6 7 "0.1328"
$newSheet.Cells.Item(6, 8).Value = 5

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: add a new top data row for 2022-Q1
#    and shift the previous rows down by one, refreshing the running index
#    in column A.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Snapshot the existing data rows (2..6) before anything gets overwritten.
# .Value2 returns a detached scalar (unlike .Value, which can behave like a
# live-bound wrapper), so this is safe to read before we start writing.
$existingRows = @()
for ($r = 2; $r -le 6; $r++) {
    $existingRows += ,@($tot.Cells.Item($r, 2).Value2, $tot.Cells.Item($r, 3).Value2, $tot.Cells.Item($r, 4).Value2)
}

# Row 7 is brand new - clone the formatting of row 6's A cell (index column
# style) onto it before writing any values.
$tot.Cells.Item(6, 1).Copy()
$tot.Cells.Item(7, 1).PasteSpecial(-4122)

# Shift the previous 5 rows down into rows 3..7.
for ($i = 0; $i -lt 5; $i++) {
    $destRow = $i + 3
    $tot.Cells.Item($destRow, 1).Value = $i + 1
    $tot.Cells.Item($destRow, 2).Value = $existingRows[$i][0]
    $tot.Cells.Item($destRow, 3).Value = $existingRows[$i][1]
    $tot.Cells.Item($destRow, 4).Value = $existingRows[$i][2]
}

# Insert the new 2022-Q1 row at the top of the data (row 2).
$tot.Cells.Item(2, 1).Value = 0
$tot.Cells.Item(2, 2).Value = "2022-Q1"
$tot.Cells.Item(2, 3).Value = 5
$tot.Cells.Item(2, 4).Value = 6.03

# Restore the originally active sheet/selection (creating/renaming sheets
# above shifts focus onto the new sheet as a side effect).
$wb.Worksheets.Item(1).Select() | Out-Null
$wb.Worksheets.Item(1).Range("A1").Select() | Out-Null
